$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above the existing "Especial/Primera/Segunda" block
# (old rows 218-223), pushing that block down to rows 221-226 and making
# room for a new weekly price block at rows 218-220.
$ws.Rows("218:220").Insert()

# New weekly block (fecha 2021-09-22 / serial 44461) at rows 218-220.
# Columns A,B,C,E,F,G,H,I,J,K,R,T are identical to the surrounding rows,
# so copy them from row 221 (the row that now holds the former row 218).
for ($r = 218; $r -le 220; $r++) {
    $ws.Cells.Item($r, 1).Value = $ws.Cells.Item(221, 1).Value()
    $ws.Cells.Item($r, 2).Value = $ws.Cells.Item(221, 2).Value()
    $ws.Cells.Item($r, 3).Value = $ws.Cells.Item(221, 3).Value()
    $ws.Cells.Item($r, 5).Value = $ws.Cells.Item(221, 5).Value()
    $ws.Cells.Item($r, 6).Value = $ws.Cells.Item(221, 6).Value()
    $ws.Cells.Item($r, 7).Value = $ws.Cells.Item(221, 7).Value()
    $ws.Cells.Item($r, 8).Value = $ws.Cells.Item(221, 8).Value()
    $ws.Cells.Item($r, 9).Value = $ws.Cells.Item(221, 9).Value()
    $ws.Cells.Item($r, 10).Value = $ws.Cells.Item(221, 10).Value()
    $ws.Cells.Item($r, 11).Value = $ws.Cells.Item(221, 11).Value()
    $ws.Cells.Item($r, 18).Value = $ws.Cells.Item(221, 18).Value()
    $ws.Cells.Item($r, 20).Value = $ws.Cells.Item(221, 20).Value()

    $ws.Cells.Item($r, 4).Value = 44461
}

# Row 218: "1a nueva(o)"
$ws.Cells.Item(218, 12).Value = "1a nueva(o)"
$ws.Cells.Item(218, 13).Value = 500
$ws.Cells.Item(218, 14).Value = 1800
$ws.Cells.Item(218, 15).Value = 1900
$ws.Cells.Item(218, 16).Value = 1850
$ws.Cells.Item(218, 17).Value = "$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(218, 19).Value = 1850

# Row 219: "2a nueva(o)"
$ws.Cells.Item(219, 12).Value = "2a nueva(o)"
$ws.Cells.Item(219, 13).Value = 400
$ws.Cells.Item(219, 14).Value = 1600
$ws.Cells.Item(219, 15).Value = 1700
$ws.Cells.Item(219, 16).Value = 1650
$ws.Cells.Item(219, 17).Value = "$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(219, 19).Value = 1650

# Row 220: "3a nueva (o)"
$ws.Cells.Item(220, 12).Value = "3a nueva (o)"
$ws.Cells.Item(220, 13).Value = 300
$ws.Cells.Item(220, 14).Value = 1300
$ws.Cells.Item(220, 15).Value = 1400
$ws.Cells.Item(220, 16).Value = 1350
$ws.Cells.Item(220, 17).Value = "$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(220, 19).Value = 1350
